$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update existing row 2 with new values ---
$ws.Range("A2").NumberFormat = "@"
$ws.Range("A2").Value = "1048018746"
$ws.Range("A2").Style = "Normal"

$ws.Range("B2").Value = "JOHAN CAMILO PEREZ SEPULVEDA"

$ws.Range("C2").NumberFormat = "@"
$ws.Range("C2").Value = "3104463513"
$ws.Range("C2").Style = "Normal"

# --- Add new row 3 ---
$ws.Range("A3").NumberFormat = "@"
$ws.Range("A3").Value = "1040328596"
$ws.Range("A3").Style = "Normal"

$ws.Range("B3").Value = "MARIA ISABEL ARANGO TOBON"

$ws.Range("C3").NumberFormat = "@"
$ws.Range("C3").Value = "3183779584"
$ws.Range("C3").Style = "Normal"

$ws.Range("D3").Value = "'"
$ws.Range("D3").Style = "Normal"

# --- Add new row 4 ---
$ws.Range("A4").NumberFormat = "@"
$ws.Range("A4").Value = "1001471005"
$ws.Range("A4").Style = "Normal"

$ws.Range("B4").Value = "MARIA LUCEIDA ZAPATA SERNA"

$ws.Range("C4").NumberFormat = "@"
$ws.Range("C4").Value = "3234676758"
$ws.Range("C4").Style = "Normal"

$ws.Range("D4").Value = "'"
$ws.Range("D4").Style = "Normal"
